$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Target: the run containing " vesces &" must become three runs:
#   " ves"  (kept formatting: color=000000, rtl=0)
#   "i"     (new formatting: rtl=0 only -- no explicit color)
#   "es &"  (same formatting as the first run: color=000000, rtl=0)
# i.e. "vesces" -> "vesies" (an "i" is inserted/substituted), with the
# inserted "i" rendered in a distinct (no-explicit-color) run.
# ------------------------------------------------------------------

# Locate a "donor" run elsewhere in the document that already carries
# the exact bare formatting we need for the inserted "i" (just
# <w:rtl w:val="0"/>, no <w:color>) -- this lets us clone that rPr
# faithfully instead of trying to (re)construct it via Font setters
# (which always serialize an explicit color value).
$donorAnchor = $d.Content
$donorAnchor.Find.Execute("propre", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$donorSearch = $d.Range($donorAnchor.End, $d.Content.End)
$donorSearch.Find.Execute("m", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$donorStart = $donorSearch.Start
$donorEnd = $donorSearch.End

# Locate the run to split.
$target = $d.Content
$target.Find.Execute(" vesces &", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$origStart = $target.Start

# Shrink the original run's text down to " ves" (drops "ces &").
$target.Text = " ves"
$splitPoint = $origStart + 4

# Insert the remaining literal text right after it.
$tail = $d.Range($splitPoint, $splitPoint)
$tail.InsertAfter("ies &")

# Re-format just the inserted "i" with the donor's (bare) run formatting.
$iRange = $d.Range($splitPoint, $splitPoint + 1)
$donor = $d.Range($donorStart, $donorEnd)
$iRange.FormattedText = $donor.FormattedText
$iRange2 = $d.Range($splitPoint, $splitPoint + 1)
$iRange2.Text = "i"

Write-Output "Result: [$($d.Range($origStart, $splitPoint + 5).Text)]"
